$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Formula = '=A23*$B$1'
$ws.Cells.Item(23, 3).Formula = '=B23/$B$4'
$ws.Cells.Item(23, 4).Formula = '=10000*C23/$B$3'

$ws.Cells.Item(24, 1).Value = 12
$ws.Cells.Item(24, 2).Formula = '=A24*$B$1'
$ws.Cells.Item(24, 3).Formula = '=B24/$B$4'
$ws.Cells.Item(24, 4).Formula = '=10000*C24/$B$3'

$ws.Cells.Item(25, 1).Value = 13
$ws.Cells.Item(25, 2).Formula = '=A25*$B$1'
$ws.Cells.Item(25, 3).Formula = '=B25/$B$4'
$ws.Cells.Item(25, 4).Formula = '=10000*C25/$B$3'

$ws.Range("F23").Select() | Out-Null
